$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '27.595.44'
    'E2' = '  -0.60%  '
    'D3' = '1.752.31'
    'E3' = '  +0.36%  '
    'D4' = '1.003'
    'E4' = '  -0.28%  '
    'D5' = '323.91'
    'E5' = '  +1.08%  '
    'D6' = '1.003'
    'E6' = '  -0.11%  '
    'D7' = '0.4595'
    'E7' = '  +8.94%  '
    'D8' = '0.3575'
    'E8' = '  -1.59%  '
    'D9' = '0.07490'
    'E9' = '  +1.48%  '
    'D10' = '42.14'
    'E10' = '  -1.99%  '
    'D11' = '1.094'
    'E11' = '  +1.13%  '
    'D12' = '1.003'
    'E12' = '  -0.23%  '
    'D13' = '20.73'
    'E13' = '  +0.50%  '
    'D14' = '6.005'
    'E14' = '  -0.46%  '
    'D15' = '7.081'
    'E15' = '  -2.14%  '
    'D16' = '1.752.52'
    'E16' = '  -1.71%  '
    'D17' = '92.28'
    'E17' = '  +1.14%  '
    'E18' = '  +1.25%  '
    'D19' = '0.06428'
    'E19' = '  +1.23%  '
    'E20' = '  -0.07%  '
    'E21' = '  -0.95%  '
    'D22' = '5.805'
    'E22' = '  -1.80%  '
    'D23' = '27.649.39'
    'E23' = '  -0.59%  '
    'D24' = '11.24'
    'E24' = '  +0.30%  '
    'D25' = '2.111'
    'E25' = '  +1.47%  '
    'D26' = '164.47'
    'E26' = '  +5.04%  '
    'D27' = '20.31'
    'E27' = '  +1.25%  '
    'D28' = '1.955.76'
    'E28' = '  -1.18%  '
    'D29' = '2.075'
    'E29' = '  -2.30%  '
    'D30' = '126.44'
    'E30' = '  +1.95%  '
    'D31' = '1.061'
    'E31' = '  -5.72%  '
    'D32' = '0.09180'
    'E32' = '  +3.64%  '
    'D33' = '3.665'
    'E33' = '  +0.61%  '
    'D34' = '5.531'
    'E34' = '  -0.20%  '
    'D35' = '11.89'
    'E35' = '  -2.61%  '
    'D36' = '0.02294'
    'E36' = '  +1.37%  '
    'D37' = '0.06050'
    'E37' = '  +1.11%  '
    'D38' = '0.2099'
    'E38' = '  +0.45%  '
    'D39' = '4.977'
    'E39' = '  +0.97%  '
    'D40' = '0.6325'
    'E40' = '  +0.80%  '
    'D41' = '1.208'
    'E41' = '  +3.33%  '
    'E42' = '  -1.13%  '
    'D43' = '7.778'
    'E43' = '  +0.00%  '
    'D44' = '13.26'
    'E44' = '  -0.68%  '
    'D45' = '0.5913'
    'E45' = '  +1.15%  '
    'D46' = '3.712'
    'E46' = '  +0.99%  '
    'D47' = '123.08'
    'E47' = '  +1.00%  '
    'D48' = '1.941'
    'E48' = '  -1.20%  '
    'E49' = '  -2.63%  '
    'E50' = '  +0.72%  '
    'D51' = '72.09'
    'E51' = '  -1.67%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = '@'
    $cell.Value = $updates[$addr]
    $cell.Style = 'Normal'
}

Write-Host "Updated $($updates.Count) cells"
